$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the Spring column course list up by one row: CPSC 4135 (currently
# standalone in A8) becomes part of the C4:C7 sequence, so each of
# CPSC 4148/4155/4157/4175 moves up one slot.
$ws.Range("C4").Value = "CPSC 4135"
$ws.Range("C5").Value = "CPSC 4148"
$ws.Range("C6").Value = "CPSC 4155"
$ws.Range("C7").Value = "CPSC 4157"

# Row 8 (old CPSC 4135 entry) is no longer needed.
$ws.Range("A8:B8").ClearContents()

# Fall 2023 course list gains CPSC 4175 at the top (row 13), pushing the
# remaining two courses down one row each.
$ws.Range("A13").Value = "CPSC 4175"
$ws.Range("B13").Value = 3
$ws.Range("A14").Value = "CPSC 4176"
$ws.Range("B14").Value = 3
$ws.Range("A15").Value = "CPSC 4000"
$ws.Range("B15").Value = 0

# The Fall/Spring/Summer 2025 section (header row 30 + totals row 38) is
# removed entirely.
$ws.Range("A30:F30").ClearContents()
$ws.Range("A38:F38").ClearContents()
